$wb = $excel.ActiveWorkbook

# --- 1. Status text update: "Ready for handoff" -> "In Translation" -----
# Touches Overview!E2, Overview!F2, zh-cn!C2, de-de!C2 (all share the same
# localization-status string).
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
}

# --- 2. Narrow the "Status" columns to fit the new, shorter text --------
# Overview columns E & F (zh-cn / de-de status) and column C ("Status") on
# the zh-cn / de-de detail sheets.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.43
$overview.Columns.Item(6).ColumnWidth = 12.43

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.43

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.43
